$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style/formatting from G1 (existing "sum" header) to the new H1 cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for the "Save" column
$ws.Range("H2").Value = 0
